$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "IDSheet"
$ws2.Range("A1").Value = "KPI ID"
$ws2.Range("A2").Value = "KC ID"
$ws2.Range("B1").Value = "aZCyzqYa2aqEtf2945cna6"
$ws2.Range("B2").Value = "524fc08b8a0e4d8d857dded88d5fb882"
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36
$null = $ws2.Range("M34").Select()

$ws1.Range("A3:A6").NumberFormat = "@"
$ws1.Range("A3").Value = "'55"
$ws1.Range("A4").Value = "'22"
$ws1.Range("A5").Value = "'56"
$ws1.Range("A6").Value = "'24"
$null = $ws1.Activate()
$null = $ws1.Range("A7").Select()
